$d = $word.ActiveDocument

# --- Change 1: split "... get out of state." into three runs, replacing
#     the final word "state" with "sync" ("... get out of sync.") ---
$oldSentence = "A reset function that could return all physical elevators to a default state, eg. ground floor, to deal with a system failure/reboot or other scenario in which the software and hardware could get out of state."

$searchRange = $d.Content
$found = $searchRange.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Re-materialize a fresh Range over the same span so InsertXML replaces
    # (rather than appends after) the matched content.
    $targetRange = $d.Range($searchRange.Start, $searchRange.End)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">A reset function that could return all physical elevators to a default state, eg. ground floor, to deal with a system failure/reboot or other scenario in which the software and hardware could get out of </w:t></w:r><w:r><w:t>sync</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $targetRange.InsertXML($xml)
}

# --- Change 2: add a new bullet to "Future Features" list after the
#     "graphical observability interface" item, about localization ---
$lastBulletText = "A graphical observability interface to represent the elevators with their states visually"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd("`r", "`a") -eq $lastBulletText) {
        $para.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.Text = "Dictionaries of phrases can be added for different languages to enable localization"
        break
    }
}
